# Replace the naive QoQ error-series values with the ifoCAST full-series
# evaluation numbers (commit: "added ifoCAST full series evaluation").
# The data block is a right-trimmed "staircase": each later row has one
# fewer trailing quarter column than the row above it. The staircase now
# starts one row earlier (row 15 loses Q9/K, row 16 loses Q8/J, ... row 24
# loses Q0/B entirely), so the previously-populated trailing cell on rows
# 15-24 is cleared in addition to every remaining cell being rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 9.06420358205617
$ws.Range("C2").Value = -8.445554111503538
$ws.Range("D2").Value = -0.7101236369231867
$ws.Range("E2").Value = 0.8845767628540272
$ws.Range("F2").Value = -2.07579547640469
$ws.Range("G2").Value = 0.1911055601680333
$ws.Range("H2").Value = -0.258641106852383
$ws.Range("I2").Value = -0.3983172652632915
$ws.Range("J2").Value = 0.1282583289796501
$ws.Range("K2").Value = -0.01340599018695743

# Row 3
$ws.Range("B3").Value = -9.361105409453614
$ws.Range("C3").Value = -1.265562551573157
$ws.Range("D3").Value = 0.5656570526562616
$ws.Range("E3").Value = -2.284997860872285
$ws.Range("F3").Value = 0.03788812835922643
$ws.Range("G3").Value = -0.3840959764584079
$ws.Range("H3").Value = -0.5098896570567474
$ws.Range("I3").Value = 0.02361057504571989
$ws.Range("J3").Value = -0.1145972687512308
$ws.Range("K3").Value = 0.5301828050662366

# Row 4
$ws.Range("B4").Value = -8.430792434160105
$ws.Range("C4").Value = -5.208997744890677
$ws.Range("D4").Value = -6.927055821804304
$ws.Range("E4").Value = -3.712639700272918
$ws.Range("F4").Value = -3.440181145272791
$ws.Range("G4").Value = -3.029594428221893
$ws.Range("H4").Value = -2.083867535901119
$ws.Range("I4").Value = -1.906303010668634
$ws.Range("J4").Value = -1.02014570283596
$ws.Range("K4").Value = -0.7105763663018719

# Row 5
$ws.Range("B5").Value = -3.528816598321453
$ws.Range("C5").Value = 0.488988940704401
$ws.Range("D5").Value = -1.464583384047825
$ws.Range("E5").Value = 0.5411798699163832
$ws.Range("F5").Value = -1.051412526311918
$ws.Range("G5").Value = 0.3281387101692014
$ws.Range("H5").Value = -0.3162926321676708
$ws.Range("I5").Value = 0.6224268434738385
$ws.Range("J5").Value = 0.5727462734054503
$ws.Range("K5").Value = 0.2775813816354755

# Row 6
$ws.Range("B6").Value = -3.238217586013594
$ws.Range("C6").Value = 0.3982872707157979
$ws.Range("D6").Value = -0.2197290542545966
$ws.Range("E6").Value = -0.501181530689666
$ws.Range("F6").Value = 0.0808854692459618
$ws.Range("G6").Value = -0.0557084915466495
$ws.Range("H6").Value = 0.5786458232441503
$ws.Range("I6").Value = 0.7044337708746514
$ws.Range("J6").Value = 0.3076113034791675
$ws.Range("K6").Value = 0.2310675656625323

# Row 7
$ws.Range("B7").Value = 0.8065197243375679
$ws.Range("C7").Value = 0.03189398963179746
$ws.Range("D7").Value = -0.5830875845128012
$ws.Range("E7").Value = 0.1559775186272474
$ws.Range("F7").Value = 0.05019866896707431
$ws.Range("G7").Value = 0.6370624694028058
$ws.Range("H7").Value = 0.770872454704743
$ws.Range("I7").Value = 0.3822932620900041
$ws.Range("J7").Value = 0.3013808575635111
$ws.Range("K7").Value = 0.577223562349573

# Row 8
$ws.Range("B8").Value = -0.07360522221701771
$ws.Range("C8").Value = -0.7091928684781591
$ws.Range("D8").Value = 0.3242255872604791
$ws.Range("E8").Value = 0.1108627615709566
$ws.Range("F8").Value = 0.6481953087774386
$ws.Range("G8").Value = 0.8337352060316735
$ws.Range("H8").Value = 0.4396447898837473
$ws.Range("I8").Value = 0.3458860345678901
$ws.Range("J8").Value = 0.6282974333759916
$ws.Range("K8").Value = 0.4593001789934875

# Row 9
$ws.Range("B9").Value = -1.313633488325898
$ws.Range("C9").Value = 0.2600222099057857
$ws.Range("D9").Value = 0.4024100696281734
$ws.Range("E9").Value = 0.610488432691505
$ws.Range("F9").Value = 0.8251511859582124
$ws.Range("G9").Value = 0.5233858674811156
$ws.Range("H9").Value = 0.3807131794959351
$ws.Range("I9").Value = 0.6548618164282887
$ws.Range("J9").Value = 0.5047990197400865
$ws.Range("K9").Value = 0.6500035257015756

# Row 10
$ws.Range("B10").Value = 0.1723813810668744
$ws.Range("C10").Value = 0.3499908675985048
$ws.Range("D10").Value = 0.6691348296765106
$ws.Range("E10").Value = 0.8300803312953138
$ws.Range("F10").Value = 0.5122925844322204
$ws.Range("G10").Value = 0.3934798787210216
$ws.Range("H10").Value = 0.6635913754756477
$ws.Range("I10").Value = 0.5073009638889382
$ws.Range("J10").Value = 0.6562794631732567
$ws.Range("K10").Value = 0.3761651213703216

# Row 11
$ws.Range("B11").Value = 0.3724188472315134
$ws.Range("C11").Value = 0.6717383692389594
$ws.Range("D11").Value = 0.8080126456412985
$ws.Range("E11").Value = 0.5063106041278653
$ws.Range("F11").Value = 0.3881366285386105
$ws.Range("G11").Value = 0.6532228421052815
$ws.Range("H11").Value = 0.4988236417653144
$ws.Range("I11").Value = 0.648507154866446
$ws.Range("J11").Value = 0.367541147837353
$ws.Range("K11").Value = 0.2069891273051463

# Row 12
$ws.Range("B12").Value = 0.6931629406961795
$ws.Range("C12").Value = 0.9281904944405108
$ws.Range("D12").Value = 0.4171162597140823
$ws.Range("E12").Value = 0.3584096246900327
$ws.Range("F12").Value = 0.6618489621351915
$ws.Range("G12").Value = 0.4728834683398792
$ws.Range("H12").Value = 0.6257285918382185
$ws.Range("I12").Value = 0.3539231126887352
$ws.Range("J12").Value = 0.1884683301199705
$ws.Range("K12").Value = 0.0637708202161355

# Row 13
$ws.Range("B13").Value = 0.8858698621201214
$ws.Range("C13").Value = 0.3879117858583491
$ws.Range("D13").Value = 0.3553890248851115
$ws.Range("E13").Value = 0.6429876785237111
$ws.Range("F13").Value = 0.4525523372746749
$ws.Range("G13").Value = 0.6111984364039706
$ws.Range("H13").Value = 0.3372820671964672
$ws.Range("I13").Value = 0.1708844846912203
$ws.Range("J13").Value = 0.04729693668609269
$ws.Range("K13").Value = 0.5643005784802194

# Row 14
$ws.Range("B14").Value = 0.7281095264762592
$ws.Range("C14").Value = 0.4307600866690486
$ws.Range("D14").Value = 0.451894392210712
$ws.Range("E14").Value = 0.4772462498614081
$ws.Range("F14").Value = 0.6205778399602666
$ws.Range("G14").Value = 0.2824712685850762
$ws.Range("H14").Value = 0.1519903938025505
$ws.Range("I14").Value = 0.03243293587104695
$ws.Range("J14").Value = 0.5355888678508487
$ws.Range("K14").Value = 0.3304428394308724

# Row 15
$ws.Range("B15").Value = 0.8814549586592991
$ws.Range("C15").Value = 0.4990050944230704
$ws.Range("D15").Value = 0.2369630190705999
$ws.Range("E15").Value = 0.6565229210261294
$ws.Range("F15").Value = 0.2804205111454833
$ws.Range("G15").Value = 0.0754221543665203
$ws.Range("H15").Value = 0.005033541495763572
$ws.Range("I15").Value = 0.5092998176531234
$ws.Range("J15").Value = 0.2870455487175398

# Row 16
$ws.Range("B16").Value = 0.8113306813466088
$ws.Range("C16").Value = 0.3702517588061242
$ws.Range("D16").Value = 0.4774444992099378
$ws.Range("E16").Value = 0.3086147499160505
$ws.Range("F16").Value = 0.1114474655589138
$ws.Range("G16").Value = -0.03250802870845004
$ws.Range("H16").Value = 0.5037574785657621
$ws.Range("I16").Value = 0.2912208776562884

# Row 17
$ws.Range("B17").Value = 0.6058422499139654
$ws.Range("C17").Value = 0.5621232400842474
$ws.Range("D17").Value = 0.1630143832823334
$ws.Range("E17").Value = 0.1214012084181564
$ws.Range("F17").Value = -0.01714285900890511
$ws.Range("G17").Value = 0.464909143608339
$ws.Range("H17").Value = 0.2753750686291025

# Row 18
$ws.Range("B18").Value = 0.8731474265656303
$ws.Range("C18").Value = 0.2801031576935286
$ws.Range("D18").Value = -0.04106293303141043
$ws.Range("E18").Value = 0.01181609145629647
$ws.Range("F18").Value = 0.5006520160632426
$ws.Range("G18").Value = 0.2440474222454754

# Row 19
$ws.Range("B19").Value = 0.526928335097643
$ws.Range("C19").Value = -0.02340706549352092
$ws.Range("D19").Value = -0.08280614765469801
$ws.Range("E19").Value = 0.5328694043117227
$ws.Range("F19").Value = 0.2564355480731927

# Row 20
$ws.Range("B20").Value = 0.216269293922691
$ws.Range("C20").Value = 0.002218095814515486
$ws.Range("D20").Value = 0.4163848259537086
$ws.Range("E20").Value = 0.2715408197250452

# Row 21
$ws.Range("B21").Value = 0.1670781728486028
$ws.Range("C21").Value = 0.4299220982016248
$ws.Range("D21").Value = 0.2085679007350822

# Row 22
$ws.Range("B22").Value = 0.6837711400670328
$ws.Range("C22").Value = 0.3078859509171186

# Row 23
$ws.Range("B23").Value = 0.3519456421565676

# Trailing cells no longer present now that the staircase shifted up a row
$ws.Range("K15").ClearContents()
$ws.Range("J16").ClearContents()
$ws.Range("I17").ClearContents()
$ws.Range("H18").ClearContents()
$ws.Range("G19").ClearContents()
$ws.Range("F20").ClearContents()
$ws.Range("E21").ClearContents()
$ws.Range("D22").ClearContents()
$ws.Range("C23").ClearContents()
$ws.Range("B24").ClearContents()
